# -----------------------------------------------------------------------------
# MAJ automatique BRVM via GitHub Actions
#
# Refreshes the BRVM market-recommendation snapshot: updated price-variation
# stats for each title/sector on "Recommandations", refreshed YTD progression
# figures on "Top_YTD", and re-sorted ranking rows 26-33 / 37-47 (a handful of
# titles moved up/down the leaderboard now that the underlying % changed), with
# a few recommendation/strategy labels flipped to match the new numbers.
# -----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations" ---
# Row 2: BRVM - SERVICES PUBLICS
$wsReco.Range("D2").Value = 3442.13
$wsReco.Range("E2").Value = 115.01
# Row 3: SAFCA CI
$wsReco.Range("D3").Value = 2760
$wsReco.Range("E3").Value = 700
# Row 4: CFAO MOTORS CI
$wsReco.Range("D4").Value = 2700
$wsReco.Range("E4").Value = 675
# Row 5: BRVM - AUTRES SECTEURS
$wsReco.Range("D5").Value = 2654.02
$wsReco.Range("E5").Value = 663.3099999999999
# Row 7: NEI-CEDA CI
$wsReco.Range("D7").Value = 2360
# Row 8: AIR LIQUIDE CI
$wsReco.Range("E8").Value = 545
# Row 9: SETAO CI
$wsReco.Range("E9").Value = 545
# Row 10: SUCRIVOIRE
$wsReco.Range("D10").Value = 1960
$wsReco.Range("E10").Value = 990
# Row 11: BRVM - TRANSPORT
$wsReco.Range("D11").Value = 1510
$wsReco.Range("E11").Value = 351.25
# Row 12: BRVM - DISTRIBUTION
$wsReco.Range("D12").Value = 1471.54
$wsReco.Range("E12").Value = 363.4
# Row 13: BRVM - AGRICULTURE
$wsReco.Range("D13").Value = 1330.06
$wsReco.Range("E13").Value = 327.55
# Row 14: BRVM - INDUSTRIE
$wsReco.Range("D14").Value = 772.11
$wsReco.Range("E14").Value = 192.14
# Row 15: BRVM-PRINCIPAL
$wsReco.Range("D15").Value = 711.92
$wsReco.Range("E15").Value = 176.98
# Row 16: BRVM - CONSOMMATION DE BASE
$wsReco.Range("D16").Value = 680.46
$wsReco.Range("E16").Value = 169.54
# Row 17: BRVM - INDUSTRIELS
$wsReco.Range("D17").Value = 568.73
$wsReco.Range("E17").Value = 133.89
# Row 18: BRVM-PRESTIGE
$wsReco.Range("D18").Value = 525.92
$wsReco.Range("E18").Value = 131.96
# Row 19: BRVM - FINANCES
$wsReco.Range("D19").Value = 495.21
$wsReco.Range("E19").Value = 123.58
# Row 20: BRVM - SERVICES FINANCIERS
$wsReco.Range("D20").Value = 486.68
$wsReco.Range("E20").Value = 121.45
# Row 21: BRVM - ENERGIE
$wsReco.Range("D21").Value = 439.06
$wsReco.Range("E21").Value = 107.92
# Row 22: BRVM - CONSOMMATION DISCRETIONNAIRE
$wsReco.Range("D22").Value = 426.64
$wsReco.Range("E22").Value = 105.86
# Row 23: BRVM - TELECOMMUNICATIONS
$wsReco.Range("D23").Value = 388.89
$wsReco.Range("E23").Value = 98.19
# Row 24: UNILEVER CI (UNLC)
$wsReco.Range("D24").Value = 22.43
$wsReco.Range("E24").Value = 7.47
# Row 26: "ECOBANK COTE D''IVOIRE (ECOC)" -> "SETAO CI (STAC)"
$wsReco.Range("A26").Value = 'SETAO CI (STAC)'
$wsReco.Range("C26").Value = 1
$wsReco.Range("D26").Value = 7.56
$wsReco.Range("E26").Value = 7.41
$wsReco.Range("G26").Value = '👀 À surveiller'
# Row 27: "SETAO CI (STAC)" -> "SAFCA CI (SAFC)"
$wsReco.Range("A27").Value = 'SAFCA CI (SAFC)'
$wsReco.Range("B27").Value = 1
$wsReco.Range("C27").Value = 0
$wsReco.Range("D27").Value = 4.55
$wsReco.Range("E27").Value = 4.55
$wsReco.Range("G27").Value = '➖ Neutre'
# Row 28: "CIE CI (CIEC)" -> "ECOBANK COTE D''IVOIRE (ECOC)"
$wsReco.Range("A28").Value = 'ECOBANK COTE D''''IVOIRE (ECOC)'
$wsReco.Range("D28").Value = 4.3
$wsReco.Range("E28").Value = 4.3
# Row 29: "SUCRIVOIRE (SCRC)" -> "CIE CI (CIEC)"
$wsReco.Range("A29").Value = 'CIE CI (CIEC)'
$wsReco.Range("D29").Value = 4
$wsReco.Range("E29").Value = 4
# Row 30: "BANK OF AFRICA BN (BOAB)" -> "SUCRIVOIRE (SCRC)"
$wsReco.Range("A30").Value = 'SUCRIVOIRE (SCRC)'
$wsReco.Range("D30").Value = 3.09
$wsReco.Range("E30").Value = 3.09
# Row 31: CFAO MOTORS CI (CFAC)
$wsReco.Range("C31").Value = 0
$wsReco.Range("D31").Value = 3.03
$wsReco.Range("G31").Value = '➖ Neutre'
# Row 32: BERNABE CI (BNBC)
$wsReco.Range("C32").Value = 0
$wsReco.Range("D32").Value = 3.02
$wsReco.Range("G32").Value = '➖ Neutre'
# Row 33: "SAFCA CI (SAFC)" -> "BANK OF AFRICA BN (BOAB)"
$wsReco.Range("A33").Value = 'BANK OF AFRICA BN (BOAB)'
$wsReco.Range("C33").Value = 0
$wsReco.Range("D33").Value = 2.81
$wsReco.Range("E33").Value = 2.81
$wsReco.Range("G33").Value = '➖ Neutre'
# Row 37: "AFRICA GLOBAL LOGISTICS CI (SDSC)" -> "ORANGE COTE D'IVOIRE (ORAC)"
$wsReco.Range("A37").Value = 'ORANGE COTE D''IVOIRE (ORAC)'
$wsReco.Range("D37").Value = -0.5600000000000001
$wsReco.Range("E37").Value = 2.6
# Row 39: "ORAGROUP TOGO (ORGT)" -> "VIVO ENERGY CI (SHEC)"
$wsReco.Range("A39").Value = 'VIVO ENERGY CI (SHEC)'
$wsReco.Range("B39").Value = 1
$wsReco.Range("C39").Value = 2
$wsReco.Range("D39").Value = -2.57
$wsReco.Range("E39").Value = 4.43
$wsReco.Range("G39").Value = '👀 À surveiller'
# Row 40: "SMB CI (SMBC)" -> "SODE CI (SDCC)"
$wsReco.Range("A40").Value = 'SODE CI (SDCC)'
$wsReco.Range("D40").Value = -2.97
$wsReco.Range("E40").Value = 2.43
# Row 41: "NEI-CEDA CI (NEIC)" -> "BANK OF AFRICA BF (BOABF)"
$wsReco.Range("A41").Value = 'BANK OF AFRICA BF (BOABF)'
$wsReco.Range("D41").Value = -5.33
$wsReco.Range("E41").Value = -5.33
# Row 42: "ORANGE COTE D'IVOIRE (ORAC)" -> "SMB CI (SMBC)"
$wsReco.Range("A42").Value = 'SMB CI (SMBC)'
$wsReco.Range("D42").Value = -6.21
$wsReco.Range("E42").Value = -6.21
# Row 43: "BANK OF AFRICA BF (BOABF)" -> "SOLIBRA CI (SLBC)"
$wsReco.Range("A43").Value = 'SOLIBRA CI (SLBC)'
$wsReco.Range("D43").Value = -6.67
$wsReco.Range("E43").Value = -6.67
# Row 44: "SODE CI (SDCC)" -> "BANK OF AFRICA NG (BOAN)"
$wsReco.Range("A44").Value = 'BANK OF AFRICA NG (BOAN)'
$wsReco.Range("D44").Value = -7.24
$wsReco.Range("E44").Value = -7.24
# Row 45: "VIVO ENERGY CI (SHEC)" -> "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$wsReco.Range("A45").Value = 'AFRICA GLOBAL LOGISTICS CI (SDSC)'
$wsReco.Range("D45").Value = -11.41
$wsReco.Range("E45").Value = -7.42
# Row 46: FILTISAC CI (FTSC)
$wsReco.Range("C46").Value = 3
$wsReco.Range("D46").Value = -16.67
$wsReco.Range("E46").Value = -7.46
$wsReco.Range("F46").Value = '🔴 Vente'
$wsReco.Range("G46").Value = '⚠️ Risque de décrochage'
# Row 47: SERVAIR ABIDJAN CI (ABJC)
$wsReco.Range("C47").Value = 3
$wsReco.Range("D47").Value = -20.09
$wsReco.Range("E47").Value = -5.17
$wsReco.Range("F47").Value = '🔴 Vente'
$wsReco.Range("G47").Value = '⚠️ Risque de décrochage'

# --- Sheet "Top_YTD" ---
# Row 2: BRVM - SERVICES PUBLICS
$wsYtd.Range("B2").Value = 10721385.23
# Row 3: SAFCA CI
$wsYtd.Range("B3").Value = 389338.4
# Row 4: CFAO MOTORS CI
$wsYtd.Range("B4").Value = 360605.15
# Row 5: BRVM - AUTRES SECTEURS
$wsYtd.Range("B5").Value = 339718.45
# Row 7: NEI-CEDA CI
$wsYtd.Range("B7").Value = 226547.41
# Row 10: BRVM - TRANSPORT
$wsYtd.Range("B10").Value = 51757.07
# Row 11: BRVM - DISTRIBUTION
$wsYtd.Range("B11").Value = 47819.64
